$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 45776.01041666666
$ws.Cells.Item(2, 2).Value = 331
$ws.Cells.Item(3, 1).Value = 45776.02083333334
$ws.Cells.Item(3, 2).Value = 341
$ws.Cells.Item(4, 1).Value = 45776.03125
$ws.Cells.Item(4, 2).Value = 333
$ws.Cells.Item(5, 1).Value = 45776.04166666666
$ws.Cells.Item(5, 2).Value = 330
$ws.Cells.Item(6, 1).Value = 45776.05208333334
$ws.Cells.Item(6, 2).Value = 237
$ws.Cells.Item(7, 1).Value = 45776.0625
$ws.Cells.Item(7, 2).Value = 235
$ws.Cells.Item(8, 1).Value = 45776.07291666666
$ws.Cells.Item(8, 2).Value = 240
$ws.Cells.Item(9, 1).Value = 45776.08333333334
$ws.Cells.Item(9, 2).Value = 228
$ws.Cells.Item(10, 1).Value = 45776.09375
$ws.Cells.Item(10, 2).Value = 176
$ws.Cells.Item(11, 1).Value = 45776.10416666666
$ws.Cells.Item(11, 2).Value = 175
$ws.Cells.Item(12, 1).Value = 45776.11458333334
$ws.Cells.Item(12, 2).Value = 179
$ws.Cells.Item(13, 1).Value = 45776.125
$ws.Cells.Item(13, 2).Value = 173
$ws.Cells.Item(14, 1).Value = 45776.13541666666
$ws.Cells.Item(14, 2).Value = 155
$ws.Cells.Item(15, 1).Value = 45776.14583333334
$ws.Cells.Item(15, 2).Value = 155
$ws.Cells.Item(16, 1).Value = 45776.15625
$ws.Cells.Item(16, 2).Value = 164
$ws.Cells.Item(17, 1).Value = 45776.16666666666
$ws.Cells.Item(17, 2).Value = 153
$ws.Cells.Item(18, 1).Value = 45776.17708333334
$ws.Cells.Item(18, 2).Value = 133
$ws.Cells.Item(19, 1).Value = 45776.1875
$ws.Cells.Item(19, 2).Value = 132
$ws.Cells.Item(20, 1).Value = 45776.19791666666
$ws.Cells.Item(20, 2).Value = 140
$ws.Cells.Item(21, 1).Value = 45776.20833333334
$ws.Cells.Item(21, 2).Value = 129
$ws.Cells.Item(22, 1).Value = 45776.21875
$ws.Cells.Item(22, 2).Value = 121
$ws.Cells.Item(23, 1).Value = 45776.22916666666
$ws.Cells.Item(23, 2).Value = 118
$ws.Cells.Item(24, 1).Value = 45776.23958333334
$ws.Cells.Item(24, 2).Value = 130
$ws.Cells.Item(25, 1).Value = 45776.25
$ws.Cells.Item(25, 2).Value = 117
$ws.Cells.Item(26, 1).Value = 45776.26041666666
$ws.Cells.Item(26, 2).Value = 96
$ws.Cells.Item(27, 1).Value = 45776.27083333334
$ws.Cells.Item(27, 2).Value = 115
$ws.Cells.Item(28, 1).Value = 45776.28125
$ws.Cells.Item(28, 2).Value = 96
$ws.Cells.Item(29, 1).Value = 45776.29166666666
$ws.Cells.Item(29, 2).Value = 100
$ws.Cells.Item(30, 1).Value = 45776.30208333334
$ws.Cells.Item(30, 2).Value = 133
$ws.Cells.Item(31, 1).Value = 45776.3125
$ws.Cells.Item(31, 2).Value = 134
$ws.Cells.Item(32, 1).Value = 45776.32291666666
$ws.Cells.Item(32, 2).Value = 138
$ws.Cells.Item(33, 1).Value = 45776.33333333334
$ws.Cells.Item(33, 2).Value = 142
$ws.Cells.Item(34, 1).Value = 45776.34375
$ws.Cells.Item(34, 2).Value = 178
$ws.Cells.Item(35, 1).Value = 45776.35416666666
$ws.Cells.Item(35, 2).Value = 184
$ws.Cells.Item(36, 1).Value = 45776.36458333334
$ws.Cells.Item(36, 2).Value = 190
$ws.Cells.Item(37, 1).Value = 45776.375
$ws.Cells.Item(37, 2).Value = 185
$ws.Cells.Item(38, 1).Value = 45776.38541666666
$ws.Cells.Item(38, 2).Value = 195
$ws.Cells.Item(39, 1).Value = 45776.39583333334
$ws.Cells.Item(39, 2).Value = 196
$ws.Cells.Item(40, 1).Value = 45776.40625
$ws.Cells.Item(40, 2).Value = 200
$ws.Cells.Item(41, 1).Value = 45776.41666666666
$ws.Cells.Item(41, 2).Value = 202
$ws.Cells.Item(42, 1).Value = 45776.42708333334
$ws.Cells.Item(42, 2).Value = 195
$ws.Cells.Item(43, 1).Value = 45776.4375
$ws.Cells.Item(43, 2).Value = 195
$ws.Cells.Item(44, 1).Value = 45776.44791666666
$ws.Cells.Item(44, 2).Value = 196
$ws.Cells.Item(45, 1).Value = 45776.45833333334
$ws.Cells.Item(45, 2).Value = 197
$ws.Cells.Item(46, 1).Value = 45776.46875
$ws.Cells.Item(46, 2).Value = 215
$ws.Cells.Item(47, 1).Value = 45776.47916666666
$ws.Cells.Item(47, 2).Value = 212
$ws.Cells.Item(48, 1).Value = 45776.48958333334
$ws.Cells.Item(48, 2).Value = 212
$ws.Cells.Item(49, 1).Value = 45776.5
$ws.Cells.Item(49, 2).Value = 212
$ws.Cells.Item(50, 1).Value = 45776.51041666666
$ws.Cells.Item(50, 2).Value = 200
$ws.Cells.Item(51, 1).Value = 45776.52083333334
$ws.Cells.Item(51, 2).Value = 197
$ws.Cells.Item(52, 1).Value = 45776.53125
$ws.Cells.Item(52, 2).Value = 197
$ws.Cells.Item(53, 1).Value = 45776.54166666666
$ws.Cells.Item(53, 2).Value = 197
$ws.Cells.Item(54, 1).Value = 45776.55208333334
$ws.Cells.Item(54, 2).Value = 207
$ws.Cells.Item(55, 1).Value = 45776.5625
$ws.Cells.Item(55, 2).Value = 206
$ws.Cells.Item(56, 1).Value = 45776.57291666666
$ws.Cells.Item(56, 2).Value = 206
$ws.Cells.Item(57, 1).Value = 45776.58333333334
$ws.Cells.Item(57, 2).Value = 206
$ws.Cells.Item(58, 1).Value = 45776.59375
$ws.Cells.Item(58, 2).Value = 186
$ws.Cells.Item(59, 1).Value = 45776.60416666666
$ws.Cells.Item(59, 2).Value = 185
$ws.Cells.Item(60, 1).Value = 45776.61458333334
$ws.Cells.Item(60, 2).Value = 185
$ws.Cells.Item(61, 1).Value = 45776.625
$ws.Cells.Item(61, 2).Value = 184
$ws.Cells.Item(62, 1).Value = 45776.63541666666
$ws.Cells.Item(62, 2).Value = 181
$ws.Cells.Item(63, 1).Value = 45776.64583333334
$ws.Cells.Item(63, 2).Value = 180
$ws.Cells.Item(64, 1).Value = 45776.65625
$ws.Cells.Item(64, 2).Value = 179
$ws.Cells.Item(65, 1).Value = 45776.66666666666
$ws.Cells.Item(65, 2).Value = 178
$ws.Cells.Item(66, 1).Value = 45776.67708333334
$ws.Cells.Item(66, 2).Value = 195
$ws.Cells.Item(67, 1).Value = 45776.6875
$ws.Cells.Item(67, 2).Value = 196
$ws.Cells.Item(68, 1).Value = 45776.69791666666
$ws.Cells.Item(68, 2).Value = 197
$ws.Cells.Item(69, 1).Value = 45776.70833333334
$ws.Cells.Item(69, 2).Value = 198
$ws.Cells.Item(70, 1).Value = 45776.71875
$ws.Cells.Item(70, 2).Value = 249
$ws.Cells.Item(71, 1).Value = 45776.72916666666
$ws.Cells.Item(71, 2).Value = 249
$ws.Cells.Item(72, 1).Value = 45776.73958333334
$ws.Cells.Item(72, 2).Value = 248
$ws.Cells.Item(73, 1).Value = 45776.75
$ws.Cells.Item(73, 2).Value = 248
$ws.Cells.Item(74, 1).Value = 45776.76041666666
$ws.Cells.Item(74, 2).Value = 318
$ws.Cells.Item(75, 1).Value = 45776.77083333334
$ws.Cells.Item(75, 2).Value = 320
$ws.Cells.Item(76, 1).Value = 45776.78125
$ws.Cells.Item(76, 2).Value = 322
$ws.Cells.Item(77, 1).Value = 45776.79166666666
$ws.Cells.Item(77, 2).Value = 324
$ws.Cells.Item(78, 1).Value = 45776.80208333334
$ws.Cells.Item(78, 2).Value = 395
$ws.Cells.Item(79, 1).Value = 45776.8125
$ws.Cells.Item(79, 2).Value = 396
$ws.Cells.Item(80, 1).Value = 45776.82291666666
$ws.Cells.Item(80, 2).Value = 399
$ws.Cells.Item(81, 1).Value = 45776.83333333334
$ws.Cells.Item(81, 2).Value = 402
$ws.Cells.Item(82, 1).Value = 45776.84375
$ws.Cells.Item(82, 2).Value = 411
$ws.Cells.Item(83, 1).Value = 45776.85416666666
$ws.Cells.Item(83, 2).Value = 414
$ws.Cells.Item(84, 1).Value = 45776.86458333334
$ws.Cells.Item(84, 2).Value = 417
$ws.Cells.Item(85, 1).Value = 45776.875
$ws.Cells.Item(85, 2).Value = 420
$ws.Cells.Item(86, 1).Value = 45776.88541666666
$ws.Cells.Item(86, 2).Value = 397
$ws.Cells.Item(87, 1).Value = 45776.89583333334
$ws.Cells.Item(87, 2).Value = 399
$ws.Cells.Item(88, 1).Value = 45776.90625
$ws.Cells.Item(88, 2).Value = 401
$ws.Cells.Item(89, 1).Value = 45776.91666666666
$ws.Cells.Item(89, 2).Value = 405
$ws.Cells.Item(90, 1).Value = 45776.92708333334
$ws.Cells.Item(90, 2).Value = 371
$ws.Cells.Item(91, 1).Value = 45776.9375
$ws.Cells.Item(91, 2).Value = 372
$ws.Cells.Item(92, 1).Value = 45776.94791666666
$ws.Cells.Item(92, 2).Value = 372
$ws.Cells.Item(93, 1).Value = 45776.95833333334
$ws.Cells.Item(93, 2).Value = 373
$ws.Cells.Item(94, 1).Value = 45776.96875
$ws.Cells.Item(94, 2).Value = 0
$ws.Cells.Item(95, 1).Value = 45776.97916666666
$ws.Cells.Item(95, 2).Value = 0
$ws.Cells.Item(96, 1).Value = 45776.98958333334
$ws.Cells.Item(96, 2).Value = 0
$ws.Cells.Item(97, 1).Value = 45777
$ws.Cells.Item(97, 2).Value = 0

Write-Host "Updated rows 2-97 with new Entsoe data"
